$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add Sunday (column H) hours for the week of row 7 (3.25 hours)
$ws.Range("H7").Value = 3.25

# Update the active cell selection to O10
$ws.Range("O10").Select()
